$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "商家编号" header (old column B) is dropped, and "商家名称" (old column C)
# is renamed to "*商家名称" and moves into column B with special styling.
$ws.Cells.Item(1, 2).Value = "*商家名称"

# Red font + centered alignment for the new B1 header.
$cellB1 = $ws.Cells.Item(1, 2)
$cellB1.Font.Color = 255
$cellB1.HorizontalAlignment = -4108  # xlCenter

# Remaining header cells in C..I keep their original text (证书类型 ... 所属机构).
$ws.Cells.Item(1, 3).Value = "证书类型"
$ws.Cells.Item(1, 4).Value = "应收单价"
$ws.Cells.Item(1, 5).Value = "数量（件）"
$ws.Cells.Item(1, 6).Value = "应收合计"
$ws.Cells.Item(1, 7).Value = "优惠金额"
$ws.Cells.Item(1, 8).Value = "实收金额"
$ws.Cells.Item(1, 9).Value = "所属机构"

# Drop the now-unused 10th column header entirely (Clear removes the cell
# record outright so the sheet dimension shrinks back to column I).
$ws.Cells.Item(1, 10).Clear()

# 所属机构 (I1) reverts to the sheet's plain default formatting (no center alignment).
$ws.Range("Z1").Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)  # xlPasteFormats

# Mirror the authored workbook's final selection state.
$ws.Range("D10").Select()
